$d = $word.ActiveDocument

# --- Split the run containing "{m" into two runs: "{" and "m" ---
$rng1 = $d.Content
$rng1.Find.Execute("{m", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng1.Find.Found) {
    # Narrow the range down to just the "m" character (the second char of the match)
    $rng1.Start = $rng1.Start + 1
    # Toggling a character-formatting property on this sub-range forces Word to break
    # it out of the enclosing run without altering its effective (inherited) formatting,
    # since we set it back to its own current/default value right away.
    $rng1.Font.Bold = $true
    $rng1.Font.Bold = $false
}

# --- Split the run containing "check()}" into two runs: "check()" and "}" ---
$rng2 = $d.Content
$rng2.Find.Execute("check()}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    # Narrow the range down to just the "}" character (the last char of the match)
    $rng2.Start = $rng2.End - 1
    $rng2.Font.Bold = $true
    $rng2.Font.Bold = $false
}
